# Update cryptos list: refreshed Price (D) / Volume(1h) (E) figures, and
# restore the correct ranking order for Cosmos / InjectiveProtocol (rows 29-30).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.807.02"
$ws.Range("E2").Value = "  -4.70%  "
$ws.Range("D3").Value = "2.325.40"
$ws.Range("E3").Value = "  -6.10%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.37"
$ws.Range("E5").Value = "  -4.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "84.44"
$ws.Range("E6").Value = "  -8.63%  "
$ws.Range("E7").Value = "  -3.58%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0811"
$ws.Range("E10").Value = "  -5.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "30.00"
$ws.Range("E11").Value = "  -8.96%  "
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "2.687.25"
$ws.Range("E13").Value = "  -5.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.40"
$ws.Range("E14").Value = "  -6.97%  "
$ws.Range("E15").Value = "  -5.43%  "
$ws.Range("D16").Value = "2.305.64"
$ws.Range("E16").Value = "  -6.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.752"
$ws.Range("E17").Value = "  -5.02%  "
$ws.Range("D18").Value = "39.795.94"
$ws.Range("E18").Value = "  -4.62%  "
$ws.Range("E19").Value = "  -4.31%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.07"
$ws.Range("E20").Value = "  -5.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "67.47"
$ws.Range("E21").Value = "  -6.04%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.60"
$ws.Range("E22").Value = "  -5.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.02"
$ws.Range("E23").Value = "  -2.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.55"
$ws.Range("E24").Value = "  -7.60%  "
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E26").Value = "  -6.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.32"
$ws.Range("E27").Value = "  -6.26%  "
$ws.Range("E28").Value = "  -1.53%  "
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.26"
$ws.Range("E29").Value = "  -4.79%  "
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.28"
$ws.Range("E30").Value = "  -2.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "152.27"
$ws.Range("E31").Value = "  -1.97%  "
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.09"
$ws.Range("E33").Value = "  -6.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.45"
$ws.Range("E34").Value = "  -4.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0721"
$ws.Range("E35").Value = "  -5.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.113"
$ws.Range("E36").Value = "  -2.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0992"
$ws.Range("E37").Value = "  -3.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.74"
$ws.Range("E38").Value = "  -6.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "15.69"
$ws.Range("E39").Value = "  -8.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.70"
$ws.Range("E40").Value = "  -7.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.81"
$ws.Range("E41").Value = "  -4.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.26"
$ws.Range("E42").Value = "  -4.92%  "
$ws.Range("D43").Value = "1.946.16"
$ws.Range("E43").Value = "  -2.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0265"
$ws.Range("E44").Value = "  -6.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "17.46"
$ws.Range("E45").Value = "  -6.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.26"
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.67"
$ws.Range("E47").Value = "  -9.83%  "
$ws.Range("D48").Value = "2.553.03"
$ws.Range("E48").Value = "  -6.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "92.64"
$ws.Range("E49").Value = "  -4.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "70.44"
$ws.Range("E50").Value = "  -7.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.30"
$ws.Range("E51").Value = "  -3.89%  "
